# GANTT A JOUR !!!!
# Update the Gantt chart: advance the scroll increment, update progress/
# durations on several milestones, add a new "QUENTIN / LOIC" assignee,
# rename a milestone, and move the current view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Activate()

# Scroll increment (E3): 9 -> 13. This is the value behind the scrollbar
# control (linked cell $E$3) and drives the recalculation of the whole
# timeline header (rows 5 & 6).
$ws.Range("E3").Value = 13

# Try to keep the scrollbar control's own cached position in sync too.
try {
    $ws.Shapes.Item("Barre de défilement 6").ControlFormat.Value = 13
} catch {
}

# Row 9 - "Rélisation et validation de la maquette de l'application"
# Avancement: 60% -> 80%
$ws.Range("D9").Value = 0.8

# Row 10 - "Création de l'arborescence"
# Nombre de jours: 16 -> 23
$ws.Range("F10").Value = 23

# Row 14 - "Création et configuration de la BDD"
# Avancement: 70% -> 100%, Nombre de jours: 1 -> 7
$ws.Range("D14").Value = 1
$ws.Range("F14").Value = 7

# Row 15 - "Création du service de connexion/déconnexion"
# Avancement: 50% -> 60%, Nombre de jours: 1 -> 7
$ws.Range("D15").Value = 0.6
$ws.Range("F15").Value = 7

# Row 18 - "Affichage planning par professeur"
# Avancement: 0% -> 10%, Début: (blank) -> 01/02/2021, Nombre de jours: 0 -> 1
$ws.Range("D18").Value = 0.1
$ws.Range("E18").Value = 44228
$ws.Range("F18").Value = 1

# Row 19 - "Gestion des classes par professeurs"
# Affecté à: QUENTIN -> QUENTIN / LOIC
# Avancement: 0% -> 10%, Début: (blank) -> 01/02/2021, Nombre de jours: 0 -> 1
$ws.Range("C19").Value = "QUENTIN / LOIC"
$ws.Range("D19").Value = 0.1
$ws.Range("E19").Value = 44228
$ws.Range("F19").Value = 1

# Row 20 - "Gestion des absences et retard des étudiants AVEC FORMULAIRE"
# Avancement cleared back to blank (was 0)
$ws.Range("D20").Value = ""

# Row 21 - "Gestion des sanctions"
# Avancement: 0% -> 10%, Début: (blank) -> 01/02/2021, Nombre de jours: 0 -> 1
$ws.Range("D21").Value = 0.1
$ws.Range("E21").Value = 44228
$ws.Range("F21").Value = 1

# Row 22 - rename milestone description
$ws.Range("B22").Value = "Rédation documentation développeur"

# Update the current view / selection shown when the workbook is reopened.
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("B23").Select()
